$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8077336666666667
$ws.Range("H2").Value = 2.423201
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.020000333333333
$ws.Range("N2").Value = 9.060001
$ws.Range("O2").Value = 0.291481777372034
$ws.Range("P2").Value = 0.291481777372034
$ws.Range("Q2").Value = 2.439355942577889
$ws.Range("R2").Value = 21.954203483201
$ws.Range("S2").Value = 0.291481777372034
$ws.Range("T2").Value = 0.291481777372034

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8077336666666667
$ws.Range("H3").Value = 2.423201
$ws.Range("O3").Value = 0.3934413518781783
$ws.Range("P3").Value = 0.3934413518781784
$ws.Range("Q3").Value = 3.292636364485111
$ws.Range("R3").Value = 29.633727280366
$ws.Range("S3").Value = 0.3934413518781783
$ws.Range("T3").Value = 0.3934413518781784

# Row 4 updates
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8077336666666667
$ws.Range("H4").Value = 2.423201
$ws.Range("M4").Value = 3.229698
$ws.Range("N4").Value = 9.689094000000001
$ws.Range("O4").Value = 0.311721195201271
$ws.Range("P4").Value = 0.3117211952012711
$ws.Range("Q4").Value = 2.608735807766001
$ws.Range("R4").Value = 23.478622269894
$ws.Range("S4").Value = 0.311721195201271
$ws.Range("T4").Value = 0.3117211952012711

# New Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tgfa"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.8077336666666667
$ws.Range("H5").Value = 2.423201
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03476766666666667
$ws.Range("N5").Value = 0.104303
$ws.Range("O5").Value = 0.003355675548516525
$ws.Range("P5").Value = 0.003355675548516525
$ws.Range("Q5").Value = 0.02808301487811111
$ws.Range("R5").Value = 0.252747133903
$ws.Range("S5").Value = 0.003355675548516525
$ws.Range("T5").Value = 0.003355675548516525
